# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 232
$wsOff.Range("C2").Value = 174
$wsOff.Range("D2").Value = 55
$wsOff.Range("E2").Value = 16
$wsOff.Range("G2").Value = 3

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 176
$wsDef.Range("C2").Value = 123
$wsDef.Range("D2").Value = 40
$wsDef.Range("E2").Value = 16
